$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new notes/cells describing the local-storage / web-table persistence flow
$ws.Range("G20").Value = "Insert Membership information into Web tables"
$ws.Range("G21").Value = "After the CC is processed, perform a qry of the cust_code to make sure it is still available"
$ws.Range("H22").Value = "if not available get new one and update web table"
$ws.Range("H23").Value = "if available continue"
$ws.Range("G24").Value = "Insert Membership information into Production tables"
$ws.Range("G25").Value = "Insert Transaction into Transaction tables with this cust_code"
$ws.Range("G26").Value = "So, half the data is being stored in web and transaction tables go straight to production"

# Update the current view/selection to reflect where the author was working
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("G27").Select()
